$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1156.2941
$ws.Range("I2").Value = 1578.7778
$ws.Range("J2").Value = 681
$ws.Range("K2").Value = 1578.7778
$ws.Range("L2").Value = 681
$ws.Range("M2").Value = -1465.7778
$ws.Range("N2").Value = -907
$ws.Range("H17").Value = 4305.7334
$ws.Range("J17").Value = 4558.2856
$ws.Range("L17").Value = 13674.8568
$ws.Range("N17").Value = -14010.8568
$ws.Range("H19").Value = 1249.7778
$ws.Range("J19").Value = 1418.6923
$ws.Range("L19").Value = 1418.6923
$ws.Range("N19").Value = -1768.6923
$ws.Range("H28").Value = 3475.125
$ws.Range("I28").Value = 2333.3333
$ws.Range("J28").Value = 4160.2
$ws.Range("K28").Value = 2333.3333
$ws.Range("L28").Value = 4160.2
$ws.Range("M28").Value = -1848.3333
$ws.Range("N28").Value = -5130.2
$ws.Range("H40").Value = 4648
$ws.Range("I40").Value = 1533.3334
$ws.Range("K40").Value = 1533.3334
$ws.Range("M40").Value = -1358.3334
$ws.Range("H43").Value = 1384.75
$ws.Range("I43").Value = 1426.5714
$ws.Range("J43").Value = 1092
$ws.Range("K43").Value = 1426.5714
$ws.Range("L43").Value = 1092
$ws.Range("M43").Value = -1357.5714
$ws.Range("N43").Value = -1230
$ws.Range("H51").Value = 7217.5
$ws.Range("I51").Value = 7766.6665
$ws.Range("J51").Value = 7139.048
$ws.Range("K51").Value = 7766.6665
$ws.Range("L51").Value = 7139.048
$ws.Range("M51").Value = -7282.6665
$ws.Range("N51").Value = -8107.048
$ws.Range("H64").Value = 6872.607
$ws.Range("I64").Value = 4968
$ws.Range("K64").Value = 4968
$ws.Range("M64").Value = -4720
$ws.Range("H67").Value = 6872.607
$ws.Range("I67").Value = 4968
$ws.Range("K67").Value = 4968
$ws.Range("M67").Value = -4110
$ws.Range("H80").Value = 657.43475
$ws.Range("I80").Value = 432.41666
$ws.Range("J80").Value = 902.9091
$ws.Range("K80").Value = 1297.24998
$ws.Range("L80").Value = 2708.7273
$ws.Range("M80").Value = -299.2499800000001
$ws.Range("N80").Value = -4704.7273
$ws.Range("H83").Value = 657.43475
$ws.Range("I83").Value = 432.41666
$ws.Range("J83").Value = 902.9091
$ws.Range("K83").Value = 3891.74994
$ws.Range("L83").Value = 8126.1819
$ws.Range("M83").Value = 1100.25006
$ws.Range("N83").Value = -18110.1819
$ws.Range("H92").Value = 1518
$ws.Range("I92").Value = 1723.6
$ws.Range("J92").Value = 832.6667
$ws.Range("K92").Value = 1723.6
$ws.Range("L92").Value = 832.6667
$ws.Range("M92").Value = -475.5999999999999
$ws.Range("N92").Value = -3328.6667
$ws.Range("H94").Value = 794
$ws.Range("I94").Value = 793.5
$ws.Range("K94").Value = 793.5
$ws.Range("M94").Value = -342.5
$ws.Range("H98").Value = 1382.5883
$ws.Range("I98").Value = 1406.5
$ws.Range("K98").Value = 1406.5
$ws.Range("M98").Value = 91.5
$ws.Range("H111").Value = 13895182
$ws.Range("I111").Value = 22230812
$ws.Range("K111").Value = 66692436
$ws.Range("M111").Value = -66689369
$ws.Range("H116").Value = 4977.5884
$ws.Range("J116").Value = 5811.4287
$ws.Range("L116").Value = 5811.4287
$ws.Range("N116").Value = -12695.4287
$ws.Range("H122").Value = 1382.5883
$ws.Range("I122").Value = 1406.5
$ws.Range("K122").Value = 4219.5
$ws.Range("M122").Value = -1769.5
$ws.Range("H127").Value = 1743
$ws.Range("I127").Value = 1281.6
$ws.Range("K127").Value = 3844.8
$ws.Range("M127").Value = 1115.2
$ws.Range("H132").Value = 31253164
$ws.Range("I132").Value = 41670130
$ws.Range("K132").Value = 125010390
$ws.Range("M132").Value = -125007860
$ws.Range("H135").Value = 973.62964
$ws.Range("I135").Value = 946.4783
$ws.Range("K135").Value = 8518.3047
$ws.Range("M135").Value = -5983.304700000001
$ws.Range("H140").Value = 119989.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 119989.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 119989.75
$ws.Range("N140").Value = -130349.75
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 1509
$ws.Range("I141").Value = 1513.0526
$ws.Range("J141").Value = 1483.3334
$ws.Range("K141").Value = 4539.1578
$ws.Range("L141").Value = 4450.0002
$ws.Range("M141").Value = 640.8422
$ws.Range("N141").Value = -14810.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 956.7273
$ws.Range("J4").Value = 538.5
$ws.Range("L4").Value = 538.5
$ws.Range("N4").Value = -770.5
$ws.Range("H32").Value = 2898.7402
$ws.Range("I32").Value = 1684.6418
$ws.Range("K32").Value = 1684.6418
$ws.Range("M32").Value = -1397.6418
$ws.Range("H45").Value = 28773104
$ws.Range("J45").Value = 3995
$ws.Range("L45").Value = 3995
$ws.Range("N45").Value = -4749
$ws.Range("H61").Value = 3730.3572
$ws.Range("I61").Value = 3180.7778
$ws.Range("J61").Value = 4719.6
$ws.Range("K61").Value = 3180.7778
$ws.Range("L61").Value = 4719.6
$ws.Range("M61").Value = -2968.7778
$ws.Range("N61").Value = -5143.6
$ws.Range("H63").Value = 2767.1667
$ws.Range("I63").Value = 2520.6
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 2520.6
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1834.6
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 2767.1667
$ws.Range("I66").Value = 2520.6
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 12603
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -9171
$ws.Range("N66").Value = -26864
$ws.Range("H97").Value = 1349299.8
$ws.Range("I97").Value = 1703975.6
$ws.Range("J97").Value = 1531
$ws.Range("K97").Value = 1703975.6
$ws.Range("L97").Value = 1531
$ws.Range("M97").Value = -1703479.6
$ws.Range("N97").Value = -2523
$ws.Range("H110").Value = 1264260.4
$ws.Range("I110").Value = 1463322.5
$ws.Range("K110").Value = 1463322.5
$ws.Range("M110").Value = -1461277.5
$ws.Range("H132").Value = 2433.2068
$ws.Range("I132").Value = 2132.238
$ws.Range("J132").Value = 3223.25
$ws.Range("K132").Value = 6396.714
$ws.Range("L132").Value = 9669.75
$ws.Range("M132").Value = -3866.714
$ws.Range("N132").Value = -14729.75
$ws.Range("H136").Value = 3730.3572
$ws.Range("I136").Value = 3180.7778
$ws.Range("J136").Value = 4719.6
$ws.Range("K136").Value = 9542.3334
$ws.Range("L136").Value = 14158.8
$ws.Range("M136").Value = -6992.3334
$ws.Range("N136").Value = -19258.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3258.3333
$ws.Range("I20").Value = 2044.5555
$ws.Range("K20").Value = 2044.5555
$ws.Range("M20").Value = -1797.5555
$ws.Range("H86").Value = 8343147
$ws.Range("I86").Value = 10011436
$ws.Range("K86").Value = 10011436
$ws.Range("M86").Value = -10010313
$ws.Range("H89").Value = 8343147
$ws.Range("I89").Value = 10011436
$ws.Range("K89").Value = 50057180
$ws.Range("M89").Value = -50051564
$ws.Range("H94").Value = 6507800
$ws.Range("I94").Value = 12988398
$ws.Range("J94").Value = 27202.143
$ws.Range("K94").Value = 12988398
$ws.Range("L94").Value = 27202.143
$ws.Range("M94").Value = -12987947
$ws.Range("N94").Value = -28104.143
$ws.Range("H105").Value = 31255000
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 2980021
$ws.Range("I107").Value = 4203794.5
$ws.Range("K107").Value = 4203794.5
$ws.Range("M107").Value = -4201874.5
$ws.Range("H134").Value = 2604.46
$ws.Range("I134").Value = 1300.6364
$ws.Range("K134").Value = 3901.9092
$ws.Range("M134").Value = -1366.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 316.2
$ws.Range("I33").Value = 316.2
$ws.Range("K33").Value = 316.2
$ws.Range("M33").Value = 62.80000000000001
$ws.Range("H58").Value = 2332.8235
$ws.Range("I58").Value = 1940.25
$ws.Range("K58").Value = 1940.25
$ws.Range("M58").Value = -1737.25
$ws.Range("H62").Value = 2409.4
$ws.Range("I62").Value = 2156.4285
$ws.Range("K62").Value = 2156.4285
$ws.Range("M62").Value = -1532.4285
$ws.Range("H65").Value = 2409.4
$ws.Range("I65").Value = 2156.4285
$ws.Range("K65").Value = 10782.1425
$ws.Range("M65").Value = -7662.1425
$ws.Range("H69").Value = 42939.555
$ws.Range("I69").Value = 14609.333
$ws.Range("J69").Value = 99600
$ws.Range("K69").Value = 14609.333
$ws.Range("L69").Value = 99600
$ws.Range("M69").Value = -13860.333
$ws.Range("N69").Value = -101098
$ws.Range("H72").Value = 42939.555
$ws.Range("I72").Value = 14609.333
$ws.Range("J72").Value = 99600
$ws.Range("K72").Value = 43827.999
$ws.Range("L72").Value = 298800
$ws.Range("M72").Value = -40083.999
$ws.Range("N72").Value = -306288
$ws.Range("H93").Value = 48179.6
$ws.Range("I93").Value = 13799.667
$ws.Range("K93").Value = 13799.667
$ws.Range("M93").Value = -11927.667
$ws.Range("H105").Value = 1935.9166
$ws.Range("I105").Value = 1374.2
$ws.Range("K105").Value = 1374.2
$ws.Range("M105").Value = 372.8
$ws.Range("H107").Value = 848.975
$ws.Range("I107").Value = 748.8823
$ws.Range("K107").Value = 748.8823
$ws.Range("M107").Value = 1171.1177
$ws.Range("H132").Value = 113584.65
$ws.Range("I132").Value = 79465.46
$ws.Range("K132").Value = 238396.38
$ws.Range("M132").Value = -235866.38
$ws.Range("H136").Value = 2332.8235
$ws.Range("I136").Value = 1940.25
$ws.Range("K136").Value = 5820.75
$ws.Range("M136").Value = -3270.75
$ws.Range("H141").Value = 50326
$ws.Range("J141").Value = 50326
$ws.Range("L141").Value = 50326
$ws.Range("N141").Value = -60686

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 52493.176
$ws.Range("J12").Value = 265.42856
$ws.Range("L12").Value = 796.28568
$ws.Range("N12").Value = -1142.28568
$ws.Range("H51").Value = 766.6667
$ws.Range("I51").Value = 320
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 960
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -500
$ws.Range("N51").Value = -9920
$ws.Range("H55").Value = 45065.434
$ws.Range("J55").Value = 171000
$ws.Range("L55").Value = 513000
$ws.Range("N55").Value = -513354
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H82").Value = 3906.5
$ws.Range("I82").Value = 2813
$ws.Range("K82").Value = 8439
$ws.Range("M82").Value = -8033
$ws.Range("H85").Value = 3906.5
$ws.Range("I85").Value = 2813
$ws.Range("K85").Value = 8439
$ws.Range("M85").Value = -7035
$ws.Range("H108").Value = 1166.7858
$ws.Range("I108").Value = 1264.3334
$ws.Range("K108").Value = 3793.0002
$ws.Range("M108").Value = -913.0002
$ws.Range("H113").Value = 2610.0322
$ws.Range("J113").Value = 1942.8334
$ws.Range("L113").Value = 5828.5002
$ws.Range("N113").Value = -10168.5002
$ws.Range("H117").Value = 2601.2307
$ws.Range("I117").Value = 1964.7142
$ws.Range("J117").Value = 3343.8333
$ws.Range("K117").Value = 5894.142599999999
$ws.Range("L117").Value = 10031.4999
$ws.Range("M117").Value = -2452.142599999999
$ws.Range("N117").Value = -16915.4999
$ws.Range("H121").Value = 5059.923
$ws.Range("I121").Value = 10238.333
$ws.Range("J121").Value = 621.2857
$ws.Range("K121").Value = 30714.999
$ws.Range("L121").Value = 1863.8571
$ws.Range("M121").Value = -29404.999
$ws.Range("N121").Value = -4483.8571
$ws.Range("H125").Value = 6872.1113
$ws.Range("I125").Value = 2800
$ws.Range("K125").Value = 8400
$ws.Range("M125").Value = -3480
$ws.Range("H129").Value = 2222783
$ws.Range("I129").Value = 2222783
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 6668349
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -6663349
$ws.Range("N129").ClearContents()
$ws.Range("H138").Value = 3399.6667
$ws.Range("I138").Value = 1999
$ws.Range("K138").Value = 5997
$ws.Range("M138").Value = -857
$ws.Range("H141").Value = 2374.5715
$ws.Range("I141").Value = 2226.6
$ws.Range("J141").Value = 2744.5
$ws.Range("K141").Value = 6679.799999999999
$ws.Range("L141").Value = 8233.5
$ws.Range("M141").Value = -1499.799999999999
$ws.Range("N141").Value = -18593.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 804.6923
$ws.Range("I2").Value = 1158.5555
$ws.Range("K2").Value = 1158.5555
$ws.Range("M2").Value = -1045.5555
$ws.Range("H21").Value = 10382.4
$ws.Range("I21").Value = 6304
$ws.Range("J21").Value = 16500
$ws.Range("K21").Value = 6304
$ws.Range("L21").Value = 16500
$ws.Range("M21").Value = -6131
$ws.Range("N21").Value = -16846
$ws.Range("H30").Value = 10382.4
$ws.Range("I30").Value = 6304
$ws.Range("J30").Value = 16500
$ws.Range("K30").Value = 6304
$ws.Range("L30").Value = 16500
$ws.Range("M30").Value = -6199
$ws.Range("N30").Value = -16710
$ws.Range("H38").Value = 27750
$ws.Range("J38").Value = 27750
$ws.Range("L38").Value = 27750
$ws.Range("N38").Value = -28676
$ws.Range("H70").Value = 20005450
$ws.Range("I70").Value = 28576214
$ws.Range("K70").Value = 28576214
$ws.Range("M70").Value = -28575944
$ws.Range("H73").Value = 20005450
$ws.Range("I73").Value = 28576214
$ws.Range("K73").Value = 28576214
$ws.Range("M73").Value = -28575278
$ws.Range("H80").Value = 1629075.5
$ws.Range("I80").Value = 4066811.8
$ws.Range("J80").Value = 3918
$ws.Range("K80").Value = 4066811.8
$ws.Range("L80").Value = 3918
$ws.Range("M80").Value = -4065813.8
$ws.Range("N80").Value = -5914
$ws.Range("H83").Value = 1629075.5
$ws.Range("I83").Value = 4066811.8
$ws.Range("J83").Value = 3918
$ws.Range("K83").Value = 20334059
$ws.Range("L83").Value = 19590
$ws.Range("M83").Value = -20329067
$ws.Range("N83").Value = -29574
$ws.Range("H102").Value = 4856869
$ws.Range("I102").Value = 6174514
$ws.Range("J102").Value = 2485107.2
$ws.Range("K102").Value = 6174514
$ws.Range("L102").Value = 2485107.2
$ws.Range("M102").Value = -6172892
$ws.Range("N102").Value = -2488351.2
$ws.Range("H119").Value = 99999
$ws.Range("J119").Value = 99999
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675
$ws.Range("H132").Value = 4084.2222
$ws.Range("I132").Value = 3531.75
$ws.Range("K132").Value = 10595.25
$ws.Range("M132").Value = -8065.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 100286.22
$ws.Range("I22").Value = 296809
$ws.Range("J22").Value = 2024.8334
$ws.Range("K22").Value = 296809
$ws.Range("L22").Value = 2024.8334
$ws.Range("M22").Value = -296514
$ws.Range("N22").Value = -2614.8334
$ws.Range("H27").Value = 100286.22
$ws.Range("I27").Value = 296809
$ws.Range("J27").Value = 2024.8334
$ws.Range("K27").Value = 296809
$ws.Range("L27").Value = 2024.8334
$ws.Range("M27").Value = -296702
$ws.Range("N27").Value = -2238.8334
$ws.Range("H33").Value = 11378.75
$ws.Range("I33").Value = 11378.75
$ws.Range("K33").Value = 11378.75
$ws.Range("M33").Value = -11088.75
$ws.Range("H55").Value = 1706.1818
$ws.Range("I55").Value = 1714.2858
$ws.Range("J55").Value = 1700.2106
$ws.Range("K55").Value = 1714.2858
$ws.Range("L55").Value = 1700.2106
$ws.Range("M55").Value = -1541.2858
$ws.Range("N55").Value = -2046.2106
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H64").Value = 78050
$ws.Range("J64").Value = 78050
$ws.Range("L64").Value = 78050
$ws.Range("N64").Value = -78500
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("H67").Value = 78050
$ws.Range("J67").Value = 78050
$ws.Range("L67").Value = 78050
$ws.Range("N67").Value = -79610
$ws.Range("H68").Value = 2532.5
$ws.Range("I68").Value = 2710
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2710
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1961
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 2532.5
$ws.Range("I71").Value = 2710
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 13550
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -9806
$ws.Range("N71").Value = -17488
$ws.Range("H93").Value = 55558520
$ws.Range("I93").Value = 111112710
$ws.Range("J93").Value = 4328.6665
$ws.Range("K93").Value = 111112710
$ws.Range("L93").Value = 4328.6665
$ws.Range("M93").Value = -111111462
$ws.Range("N93").Value = -6824.6665
$ws.Range("H132").Value = 11595
$ws.Range("I132").Value = 12654.444
$ws.Range("K132").Value = 37963.33199999999
$ws.Range("M132").Value = -35433.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7368.5957
$ws.Range("I62").Value = 3823.5293
$ws.Range("J62").Value = 9377.467
$ws.Range("K62").Value = 3823.5293
$ws.Range("L62").Value = 9377.467
$ws.Range("M62").Value = -3199.5293
$ws.Range("N62").Value = -10625.467
$ws.Range("H65").Value = 7368.5957
$ws.Range("I65").Value = 3823.5293
$ws.Range("J65").Value = 9377.467
$ws.Range("K65").Value = 19117.6465
$ws.Range("L65").Value = 46887.33500000001
$ws.Range("M65").Value = -15997.6465
$ws.Range("N65").Value = -53127.33500000001
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3000
$ws.Range("N96").Value = -5746
$ws.Range("M96").ClearContents()
$ws.Range("H106").Value = 55829.75
$ws.Range("I106").Value = 20342
$ws.Range("K106").Value = 20342
$ws.Range("M106").Value = -19080
$ws.Range("H107").Value = 111113060
$ws.Range("I107").Value = 111113060
$ws.Range("K107").Value = 333339180
$ws.Range("M107").Value = -333337260
$ws.Range("H109").Value = 59998.332
$ws.Range("J109").Value = 59998.332
$ws.Range("L109").Value = 59998.332
$ws.Range("N109").Value = -62772.332
$ws.Range("H122").Value = 1676.7307
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 3748.75
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 11246.25
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -16146.25
